$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "openweather" mini block (I26:K27) ---
# Column I written first (openweather: / json), then J26/K26, so that the
# shared-string table is populated in the same order as the target file.
$ws.Cells.Item(26, 9).Value  = "openweather:"
$ws.Cells.Item(27, 9).Value  = "json"
$ws.Cells.Item(26, 10).Value = "погода в"
$ws.Cells.Item(26, 11).Value = "Появился"

# Two timestamps (stored as Excel time serials, formatted h:mm)
$ws.Cells.Item(27, 10).Value = 0.625
$ws.Cells.Item(27, 11).Value = 0.58333333333333337
$ws.Cells.Item(28, 10).Value = 0.75
$ws.Cells.Item(28, 11).Value = 0.70833333333333337
$ws.Range("J27:K28").NumberFormat = "h:mm"

# --- Small numbered list in column I (rows 31-35) ---
$ws.Cells.Item(31, 9).Value = 1
$ws.Cells.Item(32, 9).Value = 2
$ws.Cells.Item(33, 9).Value = 3
$ws.Cells.Item(34, 9).Value = 4
$ws.Cells.Item(35, 9).Value = 5

# --- Small scratch block (row 39-40) ---
$ws.Cells.Item(39, 9).Value  = 1
$ws.Cells.Item(39, 10).Value = 2
$ws.Cells.Item(39, 11).Value = 3
$ws.Cells.Item(40, 11).Value = 2

# --- Forecast table (rows 42-48) ---
$ws.Cells.Item(42, 9).Value  = 14
$ws.Cells.Item(42, 10).Value = 1
$ws.Cells.Item(42, 11).Value = 4
$ws.Cells.Item(42, 12).Value = 2
$ws.Cells.Item(42, 13).Value = "(a-b)/2"

$ws.Cells.Item(43, 9).Value  = 15
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 3
$ws.Cells.Item(43, 12).Value = 2
$ws.Cells.Item(43, 13).Value = "(a-b)/2"

$ws.Cells.Item(44, 9).Value  = 16
$ws.Cells.Item(44, 10).Value = -1
$ws.Cells.Item(44, 11).Value = 2
$ws.Cells.Item(44, 12).Value = 2
$ws.Cells.Item(44, 13).Value = "(a-b)/2"

$ws.Cells.Item(45, 9).Value  = 17
$ws.Cells.Item(45, 10).Value = 1
$ws.Cells.Item(45, 11).Value = 4
$ws.Cells.Item(45, 12).Value = 2
$ws.Cells.Item(45, 13).Value = "(a-b)/2"

$ws.Cells.Item(46, 9).Value  = 18
$ws.Cells.Item(46, 10).Value = 0

$ws.Cells.Item(47, 9).Value  = 19
$ws.Cells.Item(47, 10).Value = -1

$ws.Cells.Item(48, 9).Value  = 20

# --- Sheet cosmetics to match the edited layout ---
# Column I grows a bit wider and loses its "best fit" auto flag.
# (Target stored width is 13.5703125 chars; the COM bridge here only keeps
# 1/6-character resolution, so 13.5 is the closest reachable value.)
$ws.Columns.Item(9).ColumnWidth = 12.666666

# Rows 2-13 (except 6,7,11) lose their explicit 14.45pt height (back to default).
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()

# Selection moves to the new entry point of the forecast table.
$ws.Range("M42").Select()
